# Updates spreadsheet with performance test results to match version 3 of
# VerySimpleXml (row 9 == "VerySimpleXml" on the Data sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Win32 columns (B=Memory, C=Load, D=Traverse, E=Query, F=Destroy) ---
$ws.Range("B9").Value = 202948
$ws.Range("C9").Value = 2030.4
$ws.Range("D9").Value = 70.150000000000006
$ws.Range("E9").Value = 5.78
$ws.Range("F9").Value = 200.68

# --- Win64 columns (H=Memory, I=Load, J=Traverse, K=Query, L=Destroy) ---
$ws.Range("H9").Value = 351184
$ws.Range("I9").Value = 1681.5
$ws.Range("J9").Value = 96.83
$ws.Range("K9").Value = 8.36
$ws.Range("L9").Value = 246.02

# The author also moved the sheet's active-cell selection down one row.
$ws.Range("A18").Select()
